# Reassess aesdd, baum1, baum2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# aesdd (row 2): update notes
$ws.Range("H2").Value = '1 "improvised" sample per speaker; actually 6 speakers'

# BAUM1 (row 4): update pos/neg counts and notes
$ws.Range("C4").Value = 227
$ws.Range("D4").Value = 620
$ws.Range("H4").Value = "contempt, surprise, and boredom mapped to negative; labels determined by interrater consensus; some of the mp4s might not have audio!; interest mapped to positive"

# BAUM2 [tr] (row 5): update pos/neg counts and notes
$ws.Range("C5").Value = 36
$ws.Range("D5").Value = 86
$ws.Range("H5").Value = "movie or TV; excluded not useful audio; labeled by majority vote; 2 excluded for label mismatch or interrater ambiguity"

# BAUM2 [en] (row 6): update pos/neg/neu counts and notes
$ws.Range("C6").Value = 49
$ws.Range("D6").Value = 13
$ws.Range("E6").Value = 13
$ws.Range("H6").Value = "movie or TV; excluded not useful audio; labeled by majority vote; 8 excluded for label mismatch or interrater ambiguity"

# EmoReact_V_1.0 (row 7): add a TODO note in a new column, highlighted yellow
$ws.Range("I7").Value = "TODO"
$ws.Range("I7").Interior.Color = 65535

# Update the active selection left by the editor
$ws.Range("F20").Select()
